# Updated cryptos list on Tue Aug 29 11:29:54 UTC 2023 with GitHub Actions
# Refresh price/volume figures (and the BitcoinCash / Avalanche row order swap)
# to match the latest Coinranking snapshot. All Coin/Link/Price/Volume(1h) cells
# hold plain text, so numeric-looking Price values are written with a "@" (Text)
# format to stop Excel from auto-converting them to numbers, then the style is
# reset back to Normal so no extra cell formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.091.77"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.648.33"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5197"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2626"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06303"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07674"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.588"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.54%  "
$ws.Range("D13").Value = "1.639.71"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "1.874.99"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5568"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "0.0₅8114"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "26.076.90"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.600"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.35%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "192.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.911"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.184"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.509"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05356"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.267"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.448"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.325"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("E34").Value = "  -2.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.417"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.779"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9413"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5595"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01572"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.779"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "1.026.34"
$ws.Range("E42").Value = "  -2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8250"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").Value = "1.785.91"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E46").Value = "  +8.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4315"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.901"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05111"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.93%  "
